$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113 (shifts existing rows 113:146 down to 114:147),
# mirroring a new weekly price entry being added into the middle of the
# chronological log kept in this sheet.
$ws.Rows(113).Insert()

$newDate = [DateTime]"2021-09-24"

$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = $newDate
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112043
$ws.Cells.Item(113, 7).Value = "Pepino ensalada"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 800
$ws.Cells.Item(113, 11).Value = 17000
$ws.Cells.Item(113, 12).Value = 19000
$ws.Cells.Item(113, 13).Value = 18000
$ws.Cells.Item(113, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 300
$ws.Cells.Item(113, 17).Value = 60
$ws.Cells.Item(113, 18).Value = "Hortaliza"
